# Update OWASP DevSecCon Summit budget with new venue information.
# Applies the commit described by the diff: revises the Assumptions sheet
# (new sponsorship/venue model, GBP conversion, t-shirts/social/website/badge
# cost lines) and rewires the "summit 2017 budget" sheet's income/expense
# rows + formulas to match, then refreshes the "as of" date.

$wb = $excel.ActiveWorkbook
$budget = $wb.Worksheets.Item("summit 2017 budget")
$assump = $wb.Worksheets.Item("Assumptions")

# ---------------------------------------------------------------------
# Assumptions sheet
# ---------------------------------------------------------------------

# Row 4: fee per attendee -> summit ticket (derived from all-in daily cost)
$assump.Range("A4").Value = "summit ticket"
$assump.Range("B4").Formula = "=B9*5"

# Row 5: venue cost per day -> owasp seed
$assump.Range("A5").Value = "owasp seed"
$assump.Range("B5").Value = 150000

# Row 6: cost per hotel room per day -> # attendees paying themselves
$assump.Range("A6").Value = "# attendees paying themselves"
$assump.Range("B6").Value = 75
$assump.Range("C6").ClearContents()

# Row 7: # attendees -> # of sponsored attendees
$assump.Range("A7").Value = "# of sponsored attendees"
$assump.Range("B7").Value = 40

# Row 8: # owasp seed -> number of summit days ; new GBP marker in D8
$assump.Range("A8").Value = "number of summit days"
$assump.Range("B8").Value = 5
$assump.Range("D8").Value = "GBP"

# Row 9: catering cost per attendee per day -> all-in cost per day per person
$assump.Range("A9").Value = "all-in cost per day per person"
$assump.Range("C9").Value = 1.24
$assump.Range("D9").Value = 250
$assump.Range("B9").Formula = "=C9*D9"

# Row 10: owasp chapters / projects(was) -> average travel per attendee
$assump.Range("A10").Value = "average travel per attendee"
$assump.Range("B10").Value = 700

# Row 11: seed fund owasp -> sponsor income
$assump.Range("A11").Value = "sponsor income"
$assump.Range("B11").Formula = "=30000"

# Row 12: -> project/chapter donactions
$assump.Range("A12").Value = "project/chapter donactions"
$assump.Range("B12").Value = 75000

# New rows 13-16
$assump.Range("A13").Value = "social event cost per attendee"
$assump.Range("B13").Value = 50

$assump.Range("A14").Value = "cost of t-shirt"
$assump.Range("B14").Value = 30

$assump.Range("A15").Value = "website cost (layout/hosting)"
$assump.Range("B15").Value = 3000

$assump.Range("A16").Value = "badge cost"
$assump.Range("B16").Value = 5

# ---------------------------------------------------------------------
# summit 2017 budget sheet
# ---------------------------------------------------------------------

# "As of" date moves out two months
$budget.Range("F2").Value = 42887

# Attendees line now sums self-pay + sponsored + volunteer counts
$budget.Range("B4").Value = "summit attendees"
$budget.Range("C4").Formula = "=Assumptions!B3+Assumptions!B6+Assumptions!B7"

# Income block
$budget.Range("B10").Value = "seed fund owasp"
$budget.Range("C10").Value = 150000

$budget.Range("B11").Value = "people covering own ticket"
$budget.Range("C11").Formula = "=Assumptions!B6*Assumptions!B8*Assumptions!B9"

$budget.Range("B12").Value = "Sponsorship"
$budget.Range("C12").Formula = "=Assumptions!B11"

$budget.Range("B13").Value = "owasp chapters / projects"
$budget.Range("C13").Formula = "=Assumptions!B12"

# Expense block
$budget.Range("E10").Value = "Venue - cost (all-in)"
$budget.Range("F10").Formula = "=(Assumptions!B6+Assumptions!B7)*Assumptions!B8*Assumptions!B9"

$budget.Range("E11").Value = "Catering"
$budget.Range("F11").Value = "included above"

$budget.Range("E12").Value = "Sponsored travel"
$budget.Range("F12").Formula = "=Assumptions!B7*Assumptions!B10"

$budget.Range("E13").ClearContents()
$budget.Range("F13").ClearContents()

$budget.Range("E14").Value = "Signage"
$budget.Range("F14").Value = 1000

$budget.Range("E15").Value = "Program"
$budget.Range("F15").Value = 1000

$budget.Range("E16").Value = "Social activities"
$budget.Range("F16").Formula = "=(Assumptions!B3+Assumptions!B6+Assumptions!B7)*Assumptions!B13"

$budget.Range("E17").Value = "t-shirts"
$budget.Range("F17").Formula = "=(Assumptions!B3+Assumptions!B7+Assumptions!B6)*Assumptions!B14"

$budget.Range("E18").Value = "Website"
$budget.Range("F18").Formula = "=Assumptions!B15"

$budget.Range("E19").Value = "Badges"
$budget.Range("F19").Formula = "=C4*Assumptions!B16"

$budget.Range("E20").Value = "Office suplies"
$budget.Range("F20").Value = 1000

$budget.Range("E21").ClearContents()
$budget.Range("E22").ClearContents()
$budget.Range("E23").ClearContents()

# Move the active selection to match the saved view
$budget.Range("F21").Select()
$assump.Range("A17").Select()
